{"js": "// Section 7 - Create Shopping Cart API\n// The author removed a leftover empty paragraph plus the trailing\n// \"Mango.Web/Views/Home/Details.cshtml\" reference paragraph near the\n// end of the document (right after the last\n// \"update-database -context ApplicationDbContext\" line).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst needle = \"Mango.Web/Views/Home/Details.cshtml\";\n\n// Locate the paragraph that contains the Mango.Web/.../Details.cshtml text.\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(needle) !== -1) {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex === -1) {\n  throw new Error(\"Could not find the 'Mango.Web/Views/Home/Details.cshtml' paragraph.\");\n}\n\n// Delete that paragraph.\nparagraphs.items[targetIndex].delete();\n\n// Delete the immediately preceding paragraph if it is empty (the blank\n// spacer paragraph that sat between the previous content line and this one).\nif (targetIndex - 1 >= 0 && paragraphs.items[targetIndex - 1].text === \"\") {\n  paragraphs.items[targetIndex - 1].delete();\n}\n\nawait context.sync();\n", "ps1": "# Section 7 - Create Shopping Cart API\n#\n# The author removed a leftover empty spacer paragraph plus the trailing\n# \"Mango.Web/Views/Home/Details.cshtml\" reference paragraph near the end\n# of the document (right after the last\n# \"update-database -context ApplicationDbContext\" line).\n\n$d = $word.ActiveDocument\n\n$needle = \"Mango.Web/Views/Home/Details.cshtml\"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = $needle\n$find.MatchWildcards = $false\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find the '$needle' paragraph.\"\n}\n\n# Expand the found hit out to its whole enclosing paragraph (this also\n# picks up the trailing paragraph mark).\n[void]$rng.Expand(4)   # wdParagraph\n\n# Build a range that addresses the paragraph immediately preceding it.\n$prevRng = $rng.Duplicate\n$prevRng.Collapse(1)        # wdCollapseStart\n[void]$prevRng.Move(1, -1)  # wdCharacter, step back one character\n[void]$prevRng.Expand(4)    # wdParagraph\n\n$prevText = $prevRng.Text.TrimEnd([char]13)\n$prevIsEmptyAndAdjacent = ($prevText -eq \"\") -and ($prevRng.End -eq $rng.Start)\n\n# Delete the later (target) paragraph first so the earlier paragraph's\n# Start/End positions stay valid for the second delete.\n$rng.Delete()\nif ($prevIsEmptyAndAdjacent) {\n    $prevRng.Delete()\n}\n"}
